$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K to E:L
$ws.Columns("D").Insert()

# Copy number formatting/styles from the (now shifted) original column E into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest period data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 279500
$ws.Range("D9").Value = 56700
$ws.Range("D10").Value = 222700
$ws.Range("D13").Value = 0
$ws.Range("D15").Value = 88800
$ws.Range("D17").Value = 150800
$ws.Range("D18").Value = 128700
$ws.Range("D20").Value = 3400
$ws.Range("D21").Value = 220900
$ws.Range("D22").Value = 50000
$ws.Range("D23").Value = 82200
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 82200
$ws.Range("D27").Value = 82200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -3400
$ws.Range("D33").Value = 82200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 82200
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 41700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 1100
$ws.Range("D44").Value = 2400
$ws.Range("D45").Value = 7100
$ws.Range("D46").Value = 52400
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1767100
$ws.Range("D49").Value = 1900
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 15500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1836800
$ws.Range("D57").Value = 4800
$ws.Range("D58").Value = 106900
$ws.Range("D59").Value = 16700
$ws.Range("D60").Value = 128400
$ws.Range("D61").Value = 970400
$ws.Range("D62").Value = 6000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1104800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 100800
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 631200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 82200
$ws.Range("D83").Value = 88800
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 148600
$ws.Range("D91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -15500
$ws.Range("D96").Value = -79300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -137400
$ws.Range("D101").Value = -200
$ws.Range("D102").Value = -4400

# "NA" text markers in new column D
$ws.Range("D12").Value = "NA"
$ws.Range("D14").Value = "NA"

# Fix a few cells whose values in the diff are not a pure column shift
$ws.Range("E27").Value = 61700
$ws.Range("E33").Value = 61700
$ws.Range("E35").Value = 61700
$ws.Range("E81").Value = 61700
